$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The schedule table (A:C) currently ends at row 27. Append the next
# entry: a date label, the period title, and the detail note.
$ws.Range("A28").Value = "10/3"
$ws.Range("B28").Value = "11/28"
$ws.Range("C28").Value = "第78期 第四代寵物"

# Mirror the author's saved cursor position on the new last row.
$ws.Range("C27").Select()
